# Apply the "Merge branch 'main' of https://github.com/usckrc/BTG" schedule update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear out the old "Education/Seurat" material & subject columns (D & E) for rows 3-9 ---
$ws.Range("D3:E9").ClearContents()

# --- Shift the "Barron and Boulpaep" reading assignments in column G down one chapter ---
$ws.Range("G3").Value = "Barron and Boulpaep Ch. 33"
$ws.Range("G4").Value = "Barron and Boulpaep Ch. 34"
$ws.Range("G5").Value = "Barron and Boulpaep Ch. 35"
$ws.Range("G6").Value = "Barron and Boulpaep Ch. 36"
$ws.Range("G7").Value = "Barron and Boulpaep Ch. 37"
$ws.Range("G8").Value = "Barron and Boulpaep Ch. 38"
$ws.Range("G9").Value = "Barron and Boulpaep Ch. 40"

# --- Add the "Ten Simple Rules" paper list under the Outline/Figures/Methods/Results block ---
# (written in the same order the strings were originally entered)
$ws.Range("F27").Value = "Ten simple rules for structuring papers.pdf"
$ws.Range("F17").Value = "Ten simple rules for reading a scientific paper.pdf"
$ws.Range("F23").Value = "Ten Simple Rules for Reproducible Computational Research.pdf"
$ws.Range("F22").Value = "Ten simple rules for biologists learning to program.pdf"
$ws.Range("F24").Value = "Ten Simple Rules for Better Figures.pdf"
$ws.Range("F21").Value = "Ten simple rules for teaching an introduction to R.pdf"
$ws.Range("F20").Value = "Ten Simple Rules for Making Good Oral Presentations.pdf"

# --- Add the new wet-lab technique list ---
$ws.Range("F29").Value = "Pippetting"
$ws.Range("F30").Value = "DNA Digest"
$ws.Range("F31").Value = "qPCR"
$ws.Range("F32").Value = "Mouse Tagging and Tailing"
$ws.Range("F33").Value = "Mouse Perfusion"
$ws.Range("F34").Value = "ELISA Albuwell"

# --- Restore the last active selection left by the editing session ---
$ws.Range("H27").Select()
